# BehaviorScenario_Info_Technology.xlsx — household.py technology list cleanup
#
# - Fix the "electric_fan" entry (drop the stray trailing space so it
#   matches the other technology keys).
# - Consolidate the separate "refrigerator" / "freezer" rows into the
#   existing "refrigerator_freezer_combi" entry: row 40 becomes
#   "refrigerator_freezer_combi" and the two now-redundant trailing rows
#   (old rows 41 "freezer" and 42 "refrigerator_freezer_combi") are removed.
# - Leave the current viewport further down the list, near the edited rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 2).Value2 = "electric_fan"
$ws.Cells.Item(40, 2).Value2 = "refrigerator_freezer_combi"
$ws.Rows("41:42").Delete()

$ws.Range("D34").Select()
